# Fruta / hortaliza, semanal
# Insert a new weekly observation at the top of the Arándano (blue) table
# (row 31), pushing the existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 31:54 down to 32:55, carrying all formatting with them.
$ws.Rows(31).Insert()

# Populate the new row 31 with the latest weekly price record.
$ws.Range("A31").Value = 4
$ws.Range("B31").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C31").Value = "Los Lagos"
$ws.Range("D31").Value = 44966
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100101
$ws.Range("H31").Value = "Berries"
$ws.Range("I31").Value = 100101001
$ws.Range("J31").Value = "Arándano (blue)"
$ws.Range("K31").Value = "Sin especificar"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 100
$ws.Range("N31").Value = 2000
$ws.Range("O31").Value = 2200
$ws.Range("P31").Value = 2100
$ws.Range("Q31").Value = "`$/bandeja 2 kilos"
$ws.Range("R31").Value = "Provincia de Curicó"
$ws.Range("S31").Value = 1050
$ws.Range("T31").Value = 2
